$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.724.32"
$ws.Range("E2").Value = "  +1.84%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.866.97"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.77"
$ws.Range("E5").Value = "  +1.55%  "
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4716"
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2754"
$ws.Range("E8").Value = "  +0.70%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06373"
$ws.Range("E9").Value = "  -0.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.72"
$ws.Range("E10").Value = "  +9.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.878.12"
$ws.Range("E11").Value = "  +1.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07455"
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.974"
$ws.Range("E13").Value = "  -0.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "85.15"
$ws.Range("E14").Value = "  +0.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6341"
$ws.Range("E15").Value = "  +0.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.690.31"
$ws.Range("E16").Value = "  +1.87%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "243.92"
$ws.Range("E17").Value = "  +5.82%  "
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.86"
$ws.Range("E19").Value = "  +1.16%  "
$ws.Range("E20").Value = "  +0.97%  "
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.014"
$ws.Range("E22").Value = "  -0.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.067"
$ws.Range("E23").Value = "  +1.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.378"
$ws.Range("E24").Value = "  +1.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "164.44"
$ws.Range("E25").Value = "  -0.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.22"
$ws.Range("E26").Value = "  +2.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.900"
$ws.Range("E27").Value = "  +1.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1018"
$ws.Range("E28").Value = "  +0.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.384"
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.082"
$ws.Range("E30").Value = "  -1.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.871"
$ws.Range("E31").Value = "  -1.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.04961"
$ws.Range("E32").Value = "  +1.48%  "
$ws.Range("E33").Value = "  +1.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7073"
$ws.Range("E34").Value = "  -1.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.714"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.01909"
$ws.Range("E36").Value = "  +1.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.689"
$ws.Range("E37").Value = "  +2.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.8846"
$ws.Range("E38").Value = "  -1.81%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.004"
$ws.Range("E39").Value = "  +2.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "105.25"
$ws.Range("E40").Value = "  -0.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.001"
$ws.Range("E41").Value = "  +0.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.588"
$ws.Range("E42").Value = "  +0.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4097"
$ws.Range("E43").Value = "  +0.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "65.80"
$ws.Range("E44").Value = "  +7.86%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.289"
$ws.Range("E45").Value = "  +3.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1220"
$ws.Range("E46").Value = "  +2.17%  "
$ws.Range("B47").Value = "Elrond"
$ws.Range("C47").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "33.78"
$ws.Range("E47").Value = "  +1.81%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.628"
$ws.Range("E48").Value = "  -0.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05573"
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.380"
$ws.Range("E50").Value = "  -1.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3702"
$ws.Range("E51").Value = "  +0.41%  "
